# Update the header row (row 1) field names to snake_case database column names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "nome_mae"
$ws.Range("C1").Value = "nome_pai"
$ws.Range("D1").Value = "data_nascimento"
$ws.Range("E1").Value = "endereco"
$ws.Range("F1").Value = "cidade"
$ws.Range("G1").Value = "tel1"
$ws.Range("H1").Value = "tel2"
$ws.Range("I1").Value = "eucaristia"
$ws.Range("J1").Value = "batismo"
$ws.Range("K1").Value = "status_crismando"
$ws.Range("L1").Value = "fk_id_catequista"

# Adjust column widths for K and L to fit the new, longer header text
$ws.Columns.Item(11).ColumnWidth = 15.88671875
$ws.Columns.Item(12).ColumnWidth = 15

# Update the view: scroll to show column C at top-left, and select L2
$ws.Range("L2").Select()
$excel.ActiveWindow.ScrollColumn = 3
